$wb = $excel.ActiveWorkbook

# --- Rename sheets: Categories -> Hierarchies, CategoriesMapping -> HierarchiesMapping ---
$wsCategories = $wb.Worksheets.Item("Categories")
$wsCategories.Name = "Hierarchies"

$wsCategoriesMapping = $wb.Worksheets.Item("CategoriesMapping")
$wsCategoriesMapping.Name = "HierarchiesMapping"

# --- Fix stray bold-duplicate style on K2/K4/K18 of "DatasetQry ds1" ---
# These cells used a style that duplicated the one already used by the
# rest of the row (e.g. J2/J4/J18); re-align them to that shared style
# so the duplicate style definition becomes unused.
$wsDs1 = $wb.Worksheets.Item("DatasetQry ds1")

$wsDs1.Range("J2").Copy()
$wsDs1.Range("K2").PasteSpecial(-4122)

$wsDs1.Range("J4").Copy()
$wsDs1.Range("K4").PasteSpecial(-4122)

$wsDs1.Range("J18").Copy()
$wsDs1.Range("K18").PasteSpecial(-4122)
